$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $ok) {
        throw "Find failed for: $oldText"
    }
    $rng.Text = $newText
}

Replace-ExactText `
    "образованием земельного участка путем раздела с сохранением в изменённых границах земельного участка с кадастровым номером 36:11:3300001:37, расположенного по адресу: Воронежская область, р-н Каменский, х Молчаново, ул Прохладная, ШРП №1" `
    "с образованием земельного участка из земель, находящихся в государственной или муниципальной собственности, расположенного по адресу: Тверская область, Старицкий район, Берновское сельское поселение, автомобильная дорога ""Берново-Воропуни""  и исправлением ошибки в местоположении границ земельного участка с кадастровым номером 69:32:0070101:401, расположенного: местоположение установлено относительно ориентира, расположенного в границах участка. Почтовый адрес ориентира: Тверская обл, р-н Старицкий, с/п Берновское, д Берново, пл Мира, д 9"

Replace-ExactText `
    "Департамент имущественных и земельных отношений Воронежской области 3666057069 1023601570904" `
    "Государственное казенное учреждение Тверской области «Дирекция территориального дорожного фонда Тверской области» 6905009018 1026900546925"

Replace-ExactText `
    "Котлярова Анна Юрьевна" `
    "Наумова Ольга Александровна"

Replace-ExactText `
    "36-13-509" `
    "36-11-185"

Replace-ExactText `
    "89066704868" `
    "8-951-540-72-04, 8(473)255-53-72"

Replace-ExactText `
    "kotl-anna@yandex.ru ООО ""Землемер"" Воронеж Фриджрха энгельса" `
    "naumovao.2011@mail.ru ООО НПП «Компьютерные технологии» 394000, г.Воронеж, ул.Ф.Энгельса, д.5"

Replace-ExactText `
    "2017-08-17" `
    "2018-04-17"
